# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that "System" is listed first among the comma-separated recorder names.
#
# Observed rule (derived from the target diff): whenever the comma-separated
# list in column G contains the exact token "System", the whole list is
# reversed (this naturally brings "System" to the front while preserving the
# relative order of the remaining names). Rows without a "System" token are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null) {
        $parts = $val -split ", "

        if ($parts -contains "System") {
            $rev = $parts[($parts.Count - 1)..0]
            $newVal = $rev -join ", "

            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
